$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the now-unused B1 cell (was styled but empty)
$ws.Range("B1").Clear()

# Move the existing entries (Don Best Sports, Dontyne Systems, Kassl)
# down three rows, to rows 5-7, without carrying any formatting.
$ws.Range("A5").Value = $ws.Range("A2").Value()
$ws.Range("A6").Value = $ws.Range("A3").Value()
$ws.Range("A7").Value = $ws.Range("A4").Value()

# Put the new company rows in the freed-up space right after the header.
$ws.Range("A2").Value = "Zoom Video Communications"
$ws.Range("A3").Value = "ZOTAC"
$ws.Range("A4").Value = "Zuken"
